# debug : [TPSCameraDirector] : forcus_point update.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 updates ---
$ws.Range("I4").Value = "11: diry_look_pos"
$ws.Range("J4").Value = 0.5
$ws.Range("L4").Value = "off"
$ws.Range("N4").Value = "Right wrist IK"

# --- Row 5 updates ---
$ws.Range("H5").Value = "RATA01"
$ws.Range("I5").Value = "11: diry_look_pos"
$ws.Range("J5").Value = 0.5
$ws.Range("L5").Value = "off"
$ws.Range("N5").Value = "Right wrist IK"

# --- Row 6 updates ---
$ws.Range("H6").Value = "RATA01"
$ws.Range("I6").Value = "11: diry_look_pos"
$ws.Range("J6").Value = 0.5
$ws.Range("N6").Value = "Right wrist IK"

# --- Row 7: new data row (right elbow IK) ---
$ws.Range("A7").Value = 2301
$ws.Range("B7").Value = "RATA04"
$ws.Range("C7").Value = 300
$ws.Range("D7").Value = "RA01"
$ws.Range("E7").Value = 1301
$ws.Range("F7").Value = "RAEE02"
$ws.Range("G7").Value = 2301
$ws.Range("H7").Value = "RATA02"
$ws.Range("I7").Value = "0: pos_to_pos"
$ws.Range("J7").Value = 0.1
$ws.Range("K7").Value = "on"
$ws.Range("L7").Value = "off"
$ws.Range("M7").Value = "on"
$ws.Range("N7").Value = "Right elbow IK"

# --- Update selection to match the author's final cursor position ---
$ws.Range("I6").Select() | Out-Null
